$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '68.045.11'
$ws.Range("E2").Value = '  +0.95%  '
$ws.Range("D3").Value = '2.627.87'
$ws.Range("E3").Value = '  -0.17%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '597.69'
$ws.Range("E5").Value = '  -0.42%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '153.56'
$ws.Range("E7").Value = '  +0.00%  '
$ws.Range("E8").Value = '  -2.04%  '
$ws.Range("D9").Value = '2.626.45'
$ws.Range("E9").Value = '  -0.18%  '
$ws.Range("E10").Value = '  +9.75%  '
$ws.Range("E11").Value = '  -0.70%  '
$ws.Range("E12").Value = '  +0.34%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.347'
$ws.Range("E13").Value = '  -1.09%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '27.64'
$ws.Range("E14").Value = '  -0.11%  '
$ws.Range("E15").Value = '  +3.66%  '
$ws.Range("D16").Value = '3.106.15'
$ws.Range("E16").Value = '  -0.22%  '
$ws.Range("D17").Value = '67.887.35'
$ws.Range("E17").Value = '  +0.73%  '
$ws.Range("D18").Value = '2.610.44'
$ws.Range("E18").Value = '  -0.73%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '374.94'
$ws.Range("E19").Value = '  +2.72%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.31'
$ws.Range("E20").Value = '  +1.05%  '
$ws.Range("E21").Value = '  -0.42%  '
$ws.Range("E22").Value = '  -1.12%  '
$ws.Range("E23").Value = '  -2.14%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.05'
$ws.Range("E24").Value = '  -3.38%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '72.52'
$ws.Range("E25").Value = '  +2.02%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.96'
$ws.Range("E27").Value = '  -2.58%  '
$ws.Range("E28").Value = '  +2.01%  '
$ws.Range("E29").Value = '  -0.26%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.00'
$ws.Range("E30").Value = '  -0.08%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '577.02'
$ws.Range("E31").Value = '  -0.88%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.40'
$ws.Range("E32").Value = '  +0.33%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '7.85'
$ws.Range("E33").Value = '  +0.50%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.84'
$ws.Range("E34").Value = '  -0.19%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.999'
$ws.Range("E35").Value = '  -0.01%  '
$ws.Range("E36").Value = '  -1.24%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.51'
$ws.Range("E37").Value = '  -1.42%  '
$ws.Range("E38").Value = '  +0.26%  '
$ws.Range("E39").Value = '  -0.68%  '
$ws.Range("E40").Value = '  +4.57%  '
$ws.Range("E41").Value = '  +0.06%  '
$ws.Range("E42").Value = '  +1.01%  '
$ws.Range("E43").Value = '  +1.40%  '
$ws.Range("D45").Value = '0.0₆0317'
$ws.Range("E45").Value = '  +10.57%  '
$ws.Range("E46").Value = '  +0.12%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '40.47'
$ws.Range("E47").Value = '  -1.88%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '155.30'
$ws.Range("E48").Value = '  -0.43%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '3.70'
$ws.Range("E49").Value = '  -0.87%  '
$ws.Range("E50").Value = '  -2.44%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '21.83'
$ws.Range("E51").Value = '  +6.57%  '
